$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the "latest net value" (最新净值) cell
$ws.Range("G2").Value = 1.3059

# Update the "set current net value" (设置当前净值) cell
$ws.Range("G11").Value = 1.2751

# Update column G width (closest achievable value to 13.75 on the
# engine's MDW-7 pixel grid)
$ws.Range("G1").ColumnWidth = 13

# Update active selection cell
$ws.Range("H8").Select()
